$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H38").Value = 3680.875
$ws.Range("I38").Value = 1999.4
$ws.Range("J38").Value = 6483.3335
$ws.Range("K38").Value = 5998.200000000001
$ws.Range("L38").Value = 19450.0005
$ws.Range("M38").Value = -5626.200000000001
$ws.Range("N38").Value = -20194.0005
$ws.Range("H53").Value = 215.5
$ws.Range("I53").Value = 234.83333
$ws.Range("K53").Value = 234.83333
$ws.Range("M53").Value = 402.16667
$ws.Range("H135").Value = 589.46155
$ws.Range("I135").Value = 589.46155
$ws.Range("K135").Value = 5305.15395
$ws.Range("M135").Value = -2770.15395
$ws.Range("H137").Value = 2760.5386
$ws.Range("I137").Value = 1867.5769
$ws.Range("J137").Value = 3653.5
$ws.Range("K137").Value = 5602.7307
$ws.Range("L137").Value = 10960.5
$ws.Range("M137").Value = -3052.7307
$ws.Range("N137").Value = -16060.5

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 77.77778000000001
$ws.Range("I5").Value = 33.666668
$ws.Range("K5").Value = 33.666668
$ws.Range("M5").Value = 78.333332
$ws.Range("H32").Value = 1637.8889
$ws.Range("I32").Value = 748.2857
$ws.Range("K32").Value = 748.2857
$ws.Range("M32").Value = -461.2857
$ws.Range("H61").Value = 3512.889
$ws.Range("I61").Value = 3512.889
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3512.889
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3300.889
$ws.Range("N61").ClearContents()
$ws.Range("H136").Value = 3512.889
$ws.Range("I136").Value = 3512.889
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 10538.667
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -7988.667000000001
$ws.Range("N136").ClearContents()

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 77.77778000000001
$ws.Range("I4").Value = 33.666668
$ws.Range("K4").Value = 33.666668
$ws.Range("M4").Value = 81.333332
$ws.Range("H19").Value = 1324.2858
$ws.Range("H81").Value = 80926.664
$ws.Range("J81").Value = 80926.664
$ws.Range("L81").Value = 80926.664
$ws.Range("N81").Value = -83048.664
$ws.Range("H84").Value = 80926.664
$ws.Range("J84").Value = 80926.664
$ws.Range("L84").Value = 242779.992
$ws.Range("N84").Value = -253387.992
$ws.Range("H86").Value = 8250
$ws.Range("I86").Value = 5666.6665
$ws.Range("J86").Value = 9800
$ws.Range("K86").Value = 5666.6665
$ws.Range("L86").Value = 9800
$ws.Range("M86").Value = -4543.6665
$ws.Range("N86").Value = -12046
$ws.Range("H87").Value = 100000
$ws.Range("J87").Value = 100000
$ws.Range("L87").Value = 100000
$ws.Range("N87").Value = -102496
$ws.Range("H89").Value = 8250
$ws.Range("I89").Value = 5666.6665
$ws.Range("J89").Value = 9800
$ws.Range("K89").Value = 28333.3325
$ws.Range("L89").Value = 49000
$ws.Range("M89").Value = -22717.3325
$ws.Range("N89").Value = -60232
$ws.Range("H90").Value = 100000
$ws.Range("J90").Value = 100000
$ws.Range("L90").Value = 300000
$ws.Range("N90").Value = -312480
$ws.Range("H105").Value = 2891.2
$ws.Range("I105").Value = 2614
$ws.Range("K105").Value = 2614
$ws.Range("M105").Value = -867
$ws.Range("H107").Value = 4717.846
$ws.Range("I107").Value = 3314.7778
$ws.Range("J107").Value = 7874.75
$ws.Range("K107").Value = 3314.7778
$ws.Range("L107").Value = 7874.75
$ws.Range("M107").Value = -1394.7778
$ws.Range("N107").Value = -11714.75

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 417.25
$ws.Range("J22").Value = 200
$ws.Range("L22").Value = 200
$ws.Range("N22").Value = -900
$ws.Range("H36").Value = 1750
$ws.Range("I36").Value = 1750
$ws.Range("K36").Value = 1750
$ws.Range("M36").Value = -1362
$ws.Range("H40").Value = 1750
$ws.Range("I40").Value = 1750
$ws.Range("K40").Value = 1750
$ws.Range("M40").Value = -1590
$ws.Range("H58").Value = 2810.1667
$ws.Range("I58").Value = 2254.2856
$ws.Range("K58").Value = 2254.2856
$ws.Range("M58").Value = -2051.2856
$ws.Range("H86").Value = 6227
$ws.Range("I86").Value = 2450
$ws.Range("J86").Value = 10004
$ws.Range("K86").Value = 2450
$ws.Range("L86").Value = 10004
$ws.Range("M86").Value = -1327
$ws.Range("N86").Value = -12250
$ws.Range("H89").Value = 6227
$ws.Range("I89").Value = 2450
$ws.Range("J89").Value = 10004
$ws.Range("K89").Value = 12250
$ws.Range("L89").Value = 50020
$ws.Range("M89").Value = -6634
$ws.Range("N89").Value = -61252
$ws.Range("H132").Value = 3416.5
$ws.Range("I132").Value = 1287.25
$ws.Range("J132").Value = 7675
$ws.Range("K132").Value = 3861.75
$ws.Range("L132").Value = 23025
$ws.Range("M132").Value = -1331.75
$ws.Range("N132").Value = -28085
$ws.Range("H136").Value = 2810.1667
$ws.Range("I136").Value = 2254.2856
$ws.Range("K136").Value = 6762.8568
$ws.Range("M136").Value = -4212.8568

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 74295.484
$ws.Range("I4").Value = 111194.22
$ws.Range("K4").Value = 333582.66
$ws.Range("M4").Value = -333470.66
$ws.Range("H12").Value = 193.54546
$ws.Range("J12").Value = 223.14285
$ws.Range("L12").Value = 669.4285500000001
$ws.Range("N12").Value = -1015.42855
$ws.Range("H139").Value = 2158.9412
$ws.Range("I139").Value = 1446.8
$ws.Range("K139").Value = 4340.4
$ws.Range("M139").Value = 799.6000000000004
$ws.Range("H140").Value = 1824.6666
$ws.Range("I140").Value = 1597.8572
$ws.Range("K140").Value = 4793.571599999999
$ws.Range("M140").Value = 386.4284000000007

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 201.90475
$ws.Range("I2").Value = 80.27273
$ws.Range("K2").Value = 80.27273
$ws.Range("M2").Value = 32.72727
$ws.Range("H11").Value = 12528526
$ws.Range("I11").Value = 12002800
$ws.Range("K11").Value = 12002800
$ws.Range("M11").Value = -12002661
$ws.Range("H122").Value = 2314.4285
$ws.Range("I122").Value = 2033.5
$ws.Range("K122").Value = 6100.5
$ws.Range("M122").Value = -3650.5
$ws.Range("H132").Value = 5621.6
$ws.Range("I132").Value = 4777.25
$ws.Range("K132").Value = 14331.75
$ws.Range("M132").Value = -11801.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1925.5714
$ws.Range("J22").Value = 2528
$ws.Range("L22").Value = 2528
$ws.Range("N22").Value = -3118
$ws.Range("H27").Value = 1925.5714
$ws.Range("J27").Value = 2528
$ws.Range("L27").Value = 2528
$ws.Range("N27").Value = -2742
$ws.Range("H132").Value = 3329.6667
$ws.Range("I132").Value = 2994.5
$ws.Range("J132").Value = 4000
$ws.Range("K132").Value = 8983.5
$ws.Range("L132").Value = 12000
$ws.Range("M132").Value = -6453.5
$ws.Range("N132").Value = -17060
$ws.Range("H139").Value = 79999.5
$ws.Range("J139").Value = 79999.5
$ws.Range("L139").Value = 79999.5
$ws.Range("N139").Value = -90279.5
